# "30 Second Script.docx" revision
#
# The author expanded two occurrences of "barcodes" into "barcode strips"
# in the second paragraph:
#
#   "...the raspberry pi, camera and barcodes. This project..."
#       -> "...the raspberry pi, camera and barcode strips. This project..."
#
#   "...Once scanned, the barcodes will be placed..."
#       -> "...Once scanned, the barcode strips will be placed..."
#
# The paragraph also still carries a "_GoBack" bookmark planted squarely in
# the middle of the word "camera" (between "cam" and "era"), a leftover from
# wherever the author's cursor last was. That bookmark must survive the
# edit, so each Find/Replace below is scoped so its match never crosses the
# bookmark - one replace finishes exactly where the bookmark starts, the
# other begins exactly where the bookmark ends.

$d = $word.ActiveDocument

# 1) Text up to (but not including) the bookmark: "...pi, cam" -> "...pi, camera and barcode"
$d.Content.Find.Execute(
    "the raspberry pi, cam", $true, $false, $false, $false, $false, $true, 1, $false,
    "the raspberry pi, camera and barcode", 2)

# 2) Text right after the bookmark: "era and barcodes." -> " strips."
$d.Content.Find.Execute(
    "era and barcodes. This project will be able to scan barcodes using the code created with",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " strips. This project will be able to scan barcodes using the code created with", 2)

# 3) Second occurrence, later in the same paragraph: "the barcodes will be place" -> "the barcode strips will be place"
$d.Content.Find.Execute(
    "Once scanned, the barcodes will be place", $true, $false, $false, $false, $false, $true, 1, $false,
    "Once scanned, the barcode strips will be place", 2)
